# Fix convert uploaded data type to google charts format. closed #9.
#
# Rewrites sheet "Data" from a 3-column numeric sample into a 4-row x
# 8-column typed data preview (name/dept/lunchTime/salary/hireDate/age/
# isSenior/seniorityStartTime), where every value except the plain numeric
# salary/age columns is stored as literal text (even though some values
# look like numbers/booleans/dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a literal text value (shared string, t="s")
# even when the text looks like a number / boolean / date / time, and do
# so without leaving a quotePrefix style behind. We compute the value via
# a formula that yields a text string ("=""...""") and then flatten the
# formula down to a plain value with Copy + PasteSpecial(xlPasteValues).
function Set-TextValue($cell, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($cell).Formula = '="' + $escaped + '"'
    $ws.Range($cell).Copy()
    $ws.Range($cell).PasteSpecial(-4163)
}

# Row 1: headers (type-annotated column names)
Set-TextValue "A1" "name(string)"
Set-TextValue "B1" "dept(string)"
Set-TextValue "C1" "lunchTime(timeofday)"
Set-TextValue "D1" "salary(number)"
Set-TextValue "E1" "hireDate(date)"
Set-TextValue "F1" "age(number)"
Set-TextValue "G1" "isSenior(boolean)"
Set-TextValue "H1" "seniorityStartTime(datetime)"

# Row 2: John
Set-TextValue "A2" "John"
Set-TextValue "B2" "Eng"
Set-TextValue "C2" "12:00:00"
$ws.Range("D2").Value = 1000
Set-TextValue "E2" "2005-03-19"
$ws.Range("F2").Value = 35
Set-TextValue "G2" "true"
Set-TextValue "H2" "2007-12-02 15:56:00"

# Row 3: Dave
Set-TextValue "A3" "Dave"
Set-TextValue "B3" "Eng"
Set-TextValue "C3" "13:01:30.123"
Set-TextValue "D3" "500.5"
Set-TextValue "E3" "2006-04-19"
$ws.Range("F3").Value = 27
Set-TextValue "G3" "false"
Set-TextValue "H3" "2005-03-09 12:30:00.32"

# Row 4: Sally
Set-TextValue "A4" "Sally"
Set-TextValue "B4" "Eng"
Set-TextValue "C4" "09:30:05"
$ws.Range("D4").Value = 600
Set-TextValue "E4" "2005-10-10"
$ws.Range("F4").Value = 30
Set-TextValue "G4" "false"
Set-TextValue "H4" "null"

Write-Host "Data sheet converted to typed google-charts preview format."
